$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Benchmarks")

# --- Gaussian filter block (rows 14/16/17, "M" = reflect-mode sample counts) ---
# These "reflect" samples are cleared out (filter method switched to 'constant').
$null = $ws.Range("M14").ClearContents()
$null = $ws.Range("M16").ClearContents()
$null = $ws.Range("M17").ClearContents()

# --- Elapsed-time formulas for the gaussian block (rows 23/25/26/28) ---
# Recomputed timings for 'constant' mode.
$ws.Range("I23").Formula = "=18.846/4"
$ws.Range("I25").Formula = "=24.748/4"
$ws.Range("I26").Formula = "=28.199/4"
$ws.Range("I28").Formula = "=26.747/4"

# --- Second gaussian benchmark block (rows 33/35/36/37) ---
$ws.Range("I33").Value = 7.758
$null = $ws.Range("M33").ClearContents()

$ws.Range("I35").Value = 10.157
$null = $ws.Range("M35").ClearContents()

$ws.Range("I36").Value = 11.526
$null = $ws.Range("M36").ClearContents()

$ws.Range("I37").Value = 10.97

# --- Selection marker left by the editor ---
$null = $ws.Range("M33:M36").Select()
